$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '24.617.00'
$ws.Range('E2').Value = '  +10.66%  '

$ws.Range('D3').Value = '1.682.48'
$ws.Range('E3').Value = '  +5.85%  '

Set-TextValue $ws 'D4' '0.9998'
$ws.Range('E4').Value = '  -0.27%  '

Set-TextValue $ws 'D5' '306.32'
$ws.Range('E5').Value = '  +2.74%  '

Set-TextValue $ws 'D6' '0.9953'
$ws.Range('E6').Value = '  +0.60%  '

Set-TextValue $ws 'D7' '0.3684'
$ws.Range('E7').Value = '  +1.81%  '

Set-TextValue $ws 'D8' '0.3411'
$ws.Range('E8').Value = '  +2.30%  '

Set-TextValue $ws 'D9' '48.34'
$ws.Range('E9').Value = '  +17.03%  '

Set-TextValue $ws 'D10' '1.161'
$ws.Range('E10').Value = '  +3.87%  '

Set-TextValue $ws 'D11' '0.07204'
$ws.Range('E11').Value = '  +3.90%  '

Set-TextValue $ws 'D12' '0.9958'
$ws.Range('E12').Value = '  -0.33%  '

Set-TextValue $ws 'D13' '6.098'
$ws.Range('E13').Value = '  +4.78%  '

Set-TextValue $ws 'D14' '20.12'
$ws.Range('E14').Value = '  +3.61%  '

Set-TextValue $ws 'D15' '6.688'
$ws.Range('E15').Value = '  +2.55%  '

$ws.Range('D16').Value = '1.679.39'
$ws.Range('E16').Value = '  +5.80%  '

$ws.Range('E17').Value = '  +3.55%  '

Set-TextValue $ws 'D18' '0.9949'
$ws.Range('E18').Value = '  +0.57%  '

Set-TextValue $ws 'D19' '0.06632'
$ws.Range('E19').Value = '  +0.83%  '

Set-TextValue $ws 'D20' '80.57'
$ws.Range('E20').Value = '  +5.96%  '

Set-TextValue $ws 'D21' '16.39'
$ws.Range('E21').Value = '  +3.59%  '

Set-TextValue $ws 'D22' '6.068'
$ws.Range('E22').Value = '  +2.76%  '

Set-TextValue $ws 'D23' '12.10'
$ws.Range('E23').Value = '  +4.42%  '

$ws.Range('D24').Value = '24.518.46'
$ws.Range('E24').Value = '  +10.38%  '

Set-TextValue $ws 'D25' '2.425'
$ws.Range('E25').Value = '  +1.95%  '

$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws 'D26' '2.652'
$ws.Range('E26').Value = '  +6.44%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D27' '152.28'
$ws.Range('E27').Value = '  +2.79%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws 'D28' '19.41'
$ws.Range('E28').Value = '  +1.52%  '

$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = '1.865.85'
$ws.Range('E29').Value = '  +6.21%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws 'D30' '127.10'
$ws.Range('E30').Value = '  +5.12%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 'D31' '6.207'
$ws.Range('E31').Value = '  +5.73%  '

$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D32' '4.020'
$ws.Range('E32').Value = '  +2.26%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D33' '0.9724'
$ws.Range('E33').Value = '  +5.98%  '

$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 'D34' '0.08400'
$ws.Range('E34').Value = '  +3.29%  '

$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws 'D35' '1.691'
$ws.Range('E35').Value = '  +4.09%  '

$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D36' '12.29'
$ws.Range('E36').Value = '  +5.18%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D37' '0.06356'
$ws.Range('E37').Value = '  +5.55%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D38' '5.293'
$ws.Range('E38').Value = '  +3.39%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D39' '8.680'
$ws.Range('E39').Value = '  +4.18%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D40' '0.02298'
$ws.Range('E40').Value = '  +5.37%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D41' '1.243'
$ws.Range('E41').Value = '  +0.48%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D42' '0.2081'
$ws.Range('E42').Value = '  +5.10%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D43' '0.6066'
$ws.Range('E43').Value = '  +5.19%  '

$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws 'D44' '0.9948'
$ws.Range('E44').Value = '  +0.51%  '

$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D45' '3.760'
$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D46' '12.97'
$ws.Range('E46').Value = '  +1.54%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws 'D47' '0.5848'
$ws.Range('E47').Value = '  +5.24%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D48' '125.07'
$ws.Range('E48').Value = '  +0.18%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D49' '1.996'
$ws.Range('E49').Value = '  +3.30%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D50' '0.07198'
$ws.Range('E50').Value = '  +7.28%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D51' '75.58'
$ws.Range('E51').Value = '  +4.62%  '
